$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.026.06"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.993.77"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'242.78"
$ws.Range("E5").Value = "  -4.84%  "
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'54.84"
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").Value = "'57.41"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "'0.0756"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").Value = "'0.0977"
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("D13").Value = "2.287.51"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "'14.12"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "'20.84"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").Value = "'0.758"
$ws.Range("E16").Value = "  -7.08%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.016.81"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'5.05"
$ws.Range("E18").Value = "  -5.25%  "
$ws.Range("D19").Value = "36.953.22"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "'68.66"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "0.0₃0810"
$ws.Range("E21").Value = "  -4.74%  "
$ws.Range("D22").Value = "'228.53"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'5.05"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("E25").Value = "  -7.31%  "
$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "'162.57"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").Value = "'19.22"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "'0.126"
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  -5.64%  "
$ws.Range("E34").Value = "  -7.62%  "
$ws.Range("D35").Value = "'4.23"
$ws.Range("E35").Value = "  -6.42%  "
$ws.Range("E36").Value = "  -5.50%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.33"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'5.30"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "'3.09"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "1.435.01"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("E44").Value = "  -5.19%  "
$ws.Range("D45").Value = "'0.0886"
$ws.Range("E45").Value = "  -8.20%  "
$ws.Range("D46").Value = "'88.37"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").Value = "'15.24"
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "'6.74"
$ws.Range("E50").Value = "  -8.05%  "
$ws.Range("D51").Value = "2.178.66"
$ws.Range("E51").Value = "  -1.87%  "
